# Auto-generated Excel COM-interop script
# Applies scheduled market-data updates to the Exodus_Profits workbook sheets
# (values refreshed for currentAveragePrice / LevePrice / LeveProfit columns)

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 86993.8
$ws.Range("J114").Value = 86993.8
$ws.Range("L114").Value = 86993.8
$ws.Range("N114").Value = -95671.8
$ws.Range("H117").Value = 95691.8
$ws.Range("J117").Value = 95691.8
$ws.Range("L117").Value = 95691.8
$ws.Range("N117").Value = -104869.8
$ws.Range("H123").Value = 86776.664
$ws.Range("J123").Value = 86776.664
$ws.Range("L123").Value = 86776.664
$ws.Range("N123").Value = -96576.664
$ws.Range("H133").Value = 92398.5
$ws.Range("J133").Value = 92398.5
$ws.Range("L133").Value = 92398.5
$ws.Range("N133").Value = -102518.5
$ws.Range("H134").Value = 60318.5
$ws.Range("J134").Value = 69791.336
$ws.Range("L134").Value = 69791.336
$ws.Range("N134").Value = -79931.336
$ws.Range("H135").Value = 684.1786
$ws.Range("J135").Value = 1522.25
$ws.Range("L135").Value = 13700.25
$ws.Range("N135").Value = -18770.25
$ws.Range("H138").Value = 17629
$ws.Range("J138").Value = 2637.0476
$ws.Range("L138").Value = 7911.1428
$ws.Range("N138").Value = -18191.1428

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 8650.5
$ws.Range("I30").Value = 11097.5
$ws.Range("J30").Value = 7834.8335
$ws.Range("K30").Value = 11097.5
$ws.Range("L30").Value = 7834.8335
$ws.Range("M30").Value = -10947.5
$ws.Range("N30").Value = -8134.8335
$ws.Range("H61").Value = 1605.1072
$ws.Range("I61").Value = 1170.9584
$ws.Range("K61").Value = 1170.9584
$ws.Range("M61").Value = -958.9584
$ws.Range("H102").Value = 23358.445
$ws.Range("I102").Value = 791.4
$ws.Range("K102").Value = 791.4
$ws.Range("M102").Value = 830.6
$ws.Range("H107").Value = 72932.8
$ws.Range("J107").Value = 72932.8
$ws.Range("L107").Value = 72932.8
$ws.Range("N107").Value = -80612.8
$ws.Range("H115").Value = 66631
$ws.Range("J115").Value = 66631
$ws.Range("L115").Value = 66631
$ws.Range("N115").Value = -69765
$ws.Range("H118").Value = 89967.836
$ws.Range("J118").Value = 89967.836
$ws.Range("L118").Value = 89967.836
$ws.Range("N118").Value = -93281.836
$ws.Range("H121").Value = 78717
$ws.Range("J121").Value = 78717
$ws.Range("L121").Value = 78717
$ws.Range("N121").Value = -82211
$ws.Range("H132").Value = 2684.9744
$ws.Range("I132").Value = 2151.1155
$ws.Range("K132").Value = 6453.3465
$ws.Range("M132").Value = -3923.3465
$ws.Range("H136").Value = 1605.1072
$ws.Range("I136").Value = 1170.9584
$ws.Range("K136").Value = 3512.8752
$ws.Range("M136").Value = -962.8751999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 54588
$ws.Range("J109").Value = 54588
$ws.Range("L109").Value = 54588
$ws.Range("N109").Value = -57362
$ws.Range("H110").Value = 63285.855
$ws.Range("J110").Value = 63285.855
$ws.Range("L110").Value = 63285.855
$ws.Range("N110").Value = -71465.85500000001
$ws.Range("H112").Value = 50084
$ws.Range("J112").Value = 52605
$ws.Range("L112").Value = 52605
$ws.Range("N112").Value = -55559
$ws.Range("H114").Value = 91930.164
$ws.Range("J114").Value = 91992
$ws.Range("L114").Value = 91992
$ws.Range("N114").Value = -100670
$ws.Range("H116").Value = 71666.664
$ws.Range("J116").Value = 71666.664
$ws.Range("L116").Value = 71666.664
$ws.Range("N116").Value = -80844.664
$ws.Range("H117").Value = 83007.57000000001
$ws.Range("J117").Value = 83007.57000000001
$ws.Range("L117").Value = 83007.57000000001
$ws.Range("N117").Value = -92185.57000000001
$ws.Range("H118").Value = 99961.39999999999
$ws.Range("J118").Value = 99961.39999999999
$ws.Range("L118").Value = 99961.39999999999
$ws.Range("N118").Value = -103275.4
$ws.Range("H119").Value = 83990.71000000001
$ws.Range("J119").Value = 83990.71000000001
$ws.Range("L119").Value = 83990.71000000001
$ws.Range("N119").Value = -93666.71000000001
$ws.Range("H122").Value = 78464.89999999999
$ws.Range("J122").Value = 78464.89999999999
$ws.Range("L122").Value = 78464.89999999999
$ws.Range("N122").Value = -88264.89999999999
$ws.Range("H127").Value = 69888
$ws.Range("J127").Value = 69888
$ws.Range("L127").Value = 69888
$ws.Range("N127").Value = -79808
$ws.Range("H132").Value = 96246
$ws.Range("J132").Value = 96246
$ws.Range("L132").Value = 96246
$ws.Range("N132").Value = -106366
$ws.Range("H134").Value = 2385.45
$ws.Range("I134").Value = 1862.9231
$ws.Range("K134").Value = 5588.7693
$ws.Range("M134").Value = -3053.7693
$ws.Range("H135").Value = 96290
$ws.Range("J135").Value = 96290
$ws.Range("L135").Value = 96290
$ws.Range("N135").Value = -106430

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 25989.857
$ws.Range("J18").Value = 25718
$ws.Range("L18").Value = 25718
$ws.Range("N18").Value = -26178
$ws.Range("H108").Value = 50310
$ws.Range("J108").Value = 50310
$ws.Range("L108").Value = 50310
$ws.Range("N108").Value = -57990
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H117").Value = 49321.5
$ws.Range("J117").Value = 49321.5
$ws.Range("L117").Value = 49321.5
$ws.Range("N117").Value = -58499.5
$ws.Range("H129").Value = 38000
$ws.Range("J129").Value = 38000
$ws.Range("L129").Value = 38000
$ws.Range("N129").Value = -48000

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1050.75
$ws.Range("I136").Value = 1050.75
$ws.Range("K136").Value = 3152.25
$ws.Range("M136").Value = 1947.75

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4222965.5
$ws.Range("J11").Value = 6000571.5
$ws.Range("L11").Value = 6000571.5
$ws.Range("N11").Value = -6000849.5
$ws.Range("H116").Value = 88562.71000000001
$ws.Range("J116").Value = 88562.71000000001
$ws.Range("L116").Value = 88562.71000000001
$ws.Range("N116").Value = -97740.71000000001
$ws.Range("H119").Value = 69460.3
$ws.Range("J119").Value = 69545.78
$ws.Range("L119").Value = 69545.78
$ws.Range("N119").Value = -79221.78

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 16068.6
$ws.Range("J97").Value = 16068.6
$ws.Range("L97").Value = 16068.6
$ws.Range("N97").Value = -18050.6
$ws.Range("H119").Value = 98884
$ws.Range("J119").Value = 98884
$ws.Range("L119").Value = 98884
$ws.Range("N119").Value = -108560
$ws.Range("H121").Value = 57832
$ws.Range("J121").Value = 63027.145
$ws.Range("L121").Value = 63027.145
$ws.Range("N121").Value = -66521.14499999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 146995.22
$ws.Range("J46").Value = 146995.22
$ws.Range("L46").Value = 146995.22
$ws.Range("N46").Value = -147457.22
$ws.Range("H121").Value = 83141.5
$ws.Range("J121").Value = 83141.5
$ws.Range("L121").Value = 83141.5
$ws.Range("N121").Value = -86635.5
$ws.Range("H132").Value = 1063085.4
$ws.Range("I132").Value = 2419.8262
$ws.Range("J132").Value = 2418380.2
$ws.Range("K132").Value = 7259.4786
$ws.Range("L132").Value = 7255140.600000001
$ws.Range("M132").Value = -4729.4786
$ws.Range("N132").Value = -7260200.600000001
$ws.Range("H134").Value = 146995.22
$ws.Range("J134").Value = 146995.22
$ws.Range("L134").Value = 440985.66
$ws.Range("N134").Value = -446055.66
$ws.Range("H136").Value = 1899.6735
$ws.Range("I136").Value = 969.6
$ws.Range("K136").Value = 2908.8
$ws.Range("M136").Value = -358.8000000000002
